$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 13 with the new hour-log entry ---
$ws.Range("A13").Value = "Added UI elements to generate maze"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 45435

$crlf = [string]([char]13) + [string]([char]10)
$descr = "Added UI elements to set already determined modifiers now through" + $crlf + `
    "buttons and sliders. These correspond to the same settings that are " + $crlf + `
    "available in the Editor. Added simple reset option that just reloads the" + $crlf + `
    "scene."
$ws.Range("D13").Value = $descr

# D13 needs wrap text enabled (matches style used by the other Description cells)
$ws.Range("D13").WrapText = $true

# Row height becomes the same "auto" wrapped height used by similar rows
$ws.Rows(13).RowHeight = 52.5

# --- Update sheet view / selection state ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("J16").Select()
